# Add two new "annulation" (goal cancellation) comment rows to REF_COMMENTS
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=id_comment, B=scenario, C=comment, D=team, E=volume, F=speed, G=language, H=speaker
$lastRow = $ws.UsedRange.Rows.Count
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = $newRow1 - 1
$ws.Cells.Item($newRow1, 2).Value = "annulation"
$ws.Cells.Item($newRow1, 3).Value = "Attention les rouges, annulation d'un but. Tout abus sera puni"
$ws.Cells.Item($newRow1, 4).Value = "rouge"
$ws.Cells.Item($newRow1, 5).Value = 0
$ws.Cells.Item($newRow1, 6).Value = 0
$ws.Cells.Item($newRow1, 7).Value = "fr_CA"
$ws.Cells.Item($newRow1, 8).Value = "Amelie"

$ws.Cells.Item($newRow2, 1).Value = $newRow2 - 1
$ws.Cells.Item($newRow2, 2).Value = "annulation"
$ws.Cells.Item($newRow2, 3).Value = "Attention les bleus, annulation d'un but. Tout abus sera puni"
$ws.Cells.Item($newRow2, 4).Value = "bleu"
$ws.Cells.Item($newRow2, 5).Value = 0
$ws.Cells.Item($newRow2, 6).Value = 0
$ws.Cells.Item($newRow2, 7).Value = "fr_CA"
$ws.Cells.Item($newRow2, 8).Value = "Amelie"

# Match the active selection recorded in the saved file (E96 selected)
$ws.Range("E96").Select()
